# Applies the "Neutron Stars" -> "Chemistry" content rewrite described by the diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Neutron Stars: Intriguing Cosmic Sentinels" "The Fascinating World of Chemistry: Unveiling the Secrets of Matter"

# Author name (was split across "Dr" + "." + " Eliza Sokolov" runs)
Replace-Text "Dr. Eliza Sokolov" "Elizabeth Palmer"

# Email local-part + domain start (was split across "eliza" + "." + "sokolov@luminaryastro" runs);
# the trailing ".org" runs are left untouched.
Replace-Text "eliza.sokolov@luminaryastro" "elizapalm2660@eduag"

# Body paragraph sentences (each already an isolated run; 1-for-1 text swap keeps formatting)
Replace-Text "Neutron stars stand alone as remarkable and perplexing celestial objects in the vast cosmic frontier" "Chemistry, the study of matter and its properties, is a captivating field that holds the key to understanding the intricate workings of our universe"

Replace-Text " Their infinitesimally minuscule sizes conceal a heart of immense density, far surpassing that of Earth's core" " This multifaceted science delves into the nature of substances, their interactions, and their transformations, unveiling a realm of knowledge that shapes our world"

Replace-Text " Formed from the collapsed remnants of massive stars, these stellar marvels pack more mass into a volume no larger than a city, defying the very essence of condensed matter" " From the air we breathe to the food we eat, chemistry plays an indispensable role in every aspect of our existence"

Replace-Text " Their surface temperatures soar to searing heights, casting them aglow with celestial radiance, while their pulsating rhythms emit cosmic symphonies that ripple through the cosmos" " Its applications span a vast spectrum, from medicine and agriculture to manufacturing and energy, impacting our lives in countless ways"

Replace-Text "In this cosmic orchestra, neutron stars serve as pivotal players, their fleeting existences governed by the interplay of gravity's unrelenting grip and the centrifugal forces born from their rapid rotation" "The study of chemistry provides us with a deeper appreciation for the intricate dance of elements, revealing the secrets of how they combine and behave"

Replace-Text " Their intense gravitational fields create a stage where time slows down and light struggles to escape, etching surreal landscapes of time dilation around these celestial wonders" " It unravels the mysteries of chemical reactions, explaining why substances undergo changes and how these changes can be harnessed for practical applications"

Replace-Text " They are gateways to fundamental physics, revealing insights into the nature of gravity, matter, and the cosmos itself, beckoning scientists to unravel their cryptic mysteries" " Through experimentation and analysis, chemists unravel the composition and properties of matter, paving the way for advancements that benefit society"

Replace-Text "The cosmos harbors various types of neutron stars, each offering a tantalizing glimpse into the kaleidoscope of celestial evolution" "Furthermore, chemistry opens doors to a world of new possibilities"

Replace-Text " Pulsars, with their rapid rotation and directional energy bursts, act like lighthouses beaconing through the cosmic darkness" " Whether synthesizing life-saving drugs, developing sustainable energy sources, or creating innovative materials, chemists are at the forefront of scientific progress, continuously pushing the boundaries of human knowledge"

# Merge of three sentences (+ their separate "." runs) into one new sentence; the final "."
# run that used to close " These stellar oddities..." is left in place to end the paragraph.
Replace-Text " Magnetars, cloaked in enigmatic magnetic fields thousand times stronger than any found on Earth, challenge our understanding of cosmic magnetism. Binary systems, where two neutron stars engage in a celestial dance, unveil gravitational interactions of cosmic proportions, probing the limits of theoretical models. These stellar oddities captivate astronomers, physicists, and astrophysicists alike, beckoning them to decipher the enigmas held within their enigmatic hearts" " By understanding the fundamental principles of chemistry, we empower ourselves to address global challenges, such as climate change, disease, and resource scarcity, and to create a more sustainable future for generations to come"

# Summary paragraph
Replace-Text "Neutron stars, with their compact and colossal mass confined within diminutive volumes, represent captivating frontiers of astrophysics" "In essence, chemistry is the science of understanding and manipulating matter, offering insights into the composition, properties, and transformations of substances"

# Merge of two runs (removes the lastRenderedPageBreak marker that used to sit mid-sentence)
Replace-Text " These stellar sentinels, forged from the cataclysmic collapse of massive stars, endure under gravity's domineering grasp and the centrifugal dance of rapid rotation" " Its applications are vast, affecting every aspect of our lives and driving scientific advancements that address global challenges"

Replace-Text " Their blistering surface temperatures and pulsed emissions make them cosmic beacons, revealing insights into fundamental physics and offering glimpses into the exquisite performances of celestial mechanics" " Chemistry empowers us to unravel the "

Replace-Text " Pulsars, magnetars, and binary neutron star systems add further intrigue to this tapestry of cosmic wonders, spurring scientific exploration into the mysteries that lie at the core of these extraordinary celestial sentinels" "mysteries of the universe, forge new materials, and develop innovative solutions to improve the human condition"

# Add a new empty trailing paragraph (matches the "+ <w:p/>" at the end of the diff)
$lastPar = $d.Paragraphs.Last
$lastPar.Range.InsertParagraphAfter()
